{"js": "// The document body is being cleared down to a single empty paragraph that\n// only keeps the hidden \"_GoBack\" bookmark (no text, no paragraph shading).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\n\n// Delete every paragraph except the very last one. The last paragraph in\n// this document is already a plain, unformatted empty paragraph, so\n// keeping it (instead of the first one, which carries shading) gives us\n// the clean target paragraph without any extra formatting clean-up.\nfor (let i = 0; i < count - 1; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Re-create the hidden \"_GoBack\" bookmark Word keeps at the last edit\n// position, anchored on the sole remaining (now empty) paragraph.\nconst remaining = body.paragraphs.getFirst();\nremaining.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document body is being cleared down to a single empty paragraph that\n# only keeps the hidden \"_GoBack\" bookmark (no text, no paragraph shading).\n$total = $d.Paragraphs.Count\n\nif ($total -gt 1) {\n    $lastPara = $d.Paragraphs($total)\n    # Deleting from the very start of the document up to (but not including)\n    # the start of the last paragraph removes all the text paragraphs and\n    # merges what remains into that last (already \"clean\") paragraph, which\n    # inherits its lack of paragraph/run formatting (no shading).\n    $rng = $d.Range(0, $lastPara.Range.Start)\n    $rng.Delete()\n}\n\n# Re-create the hidden bookmark Word keeps at the last edit position.\n$d.Bookmarks.Add(\"_GoBack\", $d.Paragraphs(1).Range)\n"}
